# Nieuwe data toegevoegd via Streamlit op 2024-12-03 18:15:02
# Appends one new inspection record (row 80) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 80

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Kleine Drift"
$ws.Cells.Item($row, 3).Value = "BSO"

# Leading apostrophe forces the report date to stay plain text (as the
# rest of the column is) instead of being auto-parsed into a date serial;
# resetting the style back to Normal afterwards drops the quote-prefix
# formatting flag so the cell keeps the sheet's default (unstyled) look.
$ws.Cells.Item($row, 4).Value = "'2024-06-19"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
